$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the "Tổng" total row (currently row 5), pushing it to row 6.
$ws.Rows.Item(5).Insert()

# --- Apply the date number format to column A (rows 2-5) ---
# First touch "yyyy-mm-dd" (lowercase) on A2 so that format gets registered (numFmtId 164),
# then switch A2 to "YYYY-MM-DD" (uppercase, numFmtId 165) which is the one actually used.
$ws.Range("A2").NumberFormat = "yyyy-mm-dd"
$ws.Range("A2").NumberFormat = "YYYY-MM-DD"
$ws.Range("A3").NumberFormat = "YYYY-MM-DD"
$ws.Range("A4").NumberFormat = "YYYY-MM-DD"
$ws.Range("A5").NumberFormat = "YYYY-MM-DD"

# --- Row 2: 2024-12-06 ---
$ws.Range("A2").Value = 45632
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 145000
$ws.Range("D2").Value = 0.123

# --- Row 3: 2024-12-19 ---
$ws.Range("A3").Value = 45645
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 180000
$ws.Range("D3").Value = 0.153

# --- Row 4: 2024-12-21 ---
$ws.Range("A4").Value = 45647
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 200000
$ws.Range("D4").Value = 0.17

# --- Row 5 (new): 2024-12-25 ---
$ws.Range("A5").Value = 45651
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 220000
$ws.Range("D5").Value = 0.187

# --- Row 6: Tổng (totals row, now shifted down from row 5) ---
$ws.Range("C6").Value = 1175000
